$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.306379666666667
$ws.Range("H2").Value = 3.919139
$ws.Range("I2").Value = 0.4034923136874173
$ws.Range("J2").Value = 0.4034923136874172
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 37.42645866666667
$ws.Range("N2").Value = 112.279376
$ws.Range("O2").Value = 0.2415534622699011
$ws.Range("P2").Value = 0.2415534622699011
$ws.Range("Q2").Value = 48.89316459747378
$ws.Range("R2").Value = 440.038481377264
$ws.Range("S2").Value = 0.09746496537048867
$ws.Range("T2").Value = 0.09746496537048864
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.306379666666667
$ws.Range("H3").Value = 3.919139
$ws.Range("I3").Value = 0.4034923136874173
$ws.Range("J3").Value = 0.4034923136874172
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 94.96115633333334
$ws.Range("N3").Value = 284.883469
$ws.Range("O3").Value = 0.6128871635375853
$ws.Range("P3").Value = 0.6128871635375853
$ws.Range("Q3").Value = 124.0553237570212
$ws.Range("R3").Value = 1116.497913813191
$ws.Range("S3").Value = 0.2472952596450988
$ws.Range("T3").Value = 0.2472952596450988
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.306379666666667
$ws.Range("H4").Value = 3.919139
$ws.Range("I4").Value = 0.4034923136874173
$ws.Range("J4").Value = 0.4034923136874172
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 22.553069
$ws.Range("N4").Value = 67.65920699999999
$ws.Range("O4").Value = 0.1455593741925136
$ws.Range("P4").Value = 0.1455593741925136
$ws.Range("Q4").Value = 29.46287076253033
$ws.Range("R4").Value = 265.1658368627729
$ws.Range("S4").Value = 0.05873208867182984
$ws.Range("T4").Value = 0.05873208867182982
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.072562
$ws.Range("H5").Value = 3.217686
$ws.Range("I5").Value = 0.3312746929515923
$ws.Range("J5").Value = 0.3312746929515923
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 37.42645866666667
$ws.Range("N5").Value = 112.279376
$ws.Range("O5").Value = 0.2415534622699011
$ws.Range("P5").Value = 0.2415534622699011
$ws.Range("Q5").Value = 40.14219736043733
$ws.Range("R5").Value = 361.2797762439359
$ws.Range("S5").Value = 0.08002054904485553
$ws.Range("T5").Value = 0.08002054904485553
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.072562
$ws.Range("H6").Value = 3.217686
$ws.Range("I6").Value = 0.3312746929515923
$ws.Range("J6").Value = 0.3312746929515923
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 94.96115633333334
$ws.Range("N6").Value = 284.883469
$ws.Range("O6").Value = 0.6128871635375853
$ws.Range("P6").Value = 0.6128871635375853
$ws.Range("Q6").Value = 101.8517277591927
$ws.Range("R6").Value = 916.6655498327339
$ws.Range("S6").Value = 0.2030340069148859
$ws.Range("T6").Value = 0.2030340069148859
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.072562
$ws.Range("H7").Value = 3.217686
$ws.Range("I7").Value = 0.3312746929515923
$ws.Range("J7").Value = 0.3312746929515923
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 22.553069
$ws.Range("N7").Value = 67.65920699999999
$ws.Range("O7").Value = 0.1455593741925136
$ws.Range("P7").Value = 0.1455593741925136
$ws.Range("Q7").Value = 24.18956479277799
$ws.Range("R7").Value = 217.706083135002
$ws.Range("S7").Value = 0.04822013699185086
$ws.Range("T7").Value = 0.04822013699185086
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.8587400000000001
$ws.Range("H8").Value = 2.57622
$ws.Range("I8").Value = 0.2652329933609903
$ws.Range("J8").Value = 0.2652329933609903
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 37.42645866666667
$ws.Range("N8").Value = 112.279376
$ws.Range("O8").Value = 0.2415534622699011
$ws.Range("P8").Value = 0.2415534622699011
$ws.Range("Q8").Value = 32.13959711541334
$ws.Range("R8").Value = 289.25637403872
$ws.Range("S8").Value = 0.0640679478545569
$ws.Range("T8").Value = 0.06406794785455688
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.8587400000000001
$ws.Range("H9").Value = 2.57622
$ws.Range("I9").Value = 0.2652329933609903
$ws.Range("J9").Value = 0.2652329933609903
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 94.96115633333334
$ws.Range("N9").Value = 284.883469
$ws.Range("O9").Value = 0.6128871635375853
$ws.Range("P9").Value = 0.6128871635375853
$ws.Range("Q9").Value = 81.54694338968667
$ws.Range("R9").Value = 733.92249050718
$ws.Range("S9").Value = 0.1625578969776006
$ws.Range("T9").Value = 0.1625578969776005
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.8587400000000001
$ws.Range("H10").Value = 2.57622
$ws.Range("I10").Value = 0.2652329933609903
$ws.Range("J10").Value = 0.2652329933609903
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 22.553069
$ws.Range("N10").Value = 67.65920699999999
$ws.Range("O10").Value = 0.1455593741925136
$ws.Range("P10").Value = 0.1455593741925136
$ws.Range("Q10").Value = 19.36722247306
$ws.Range("R10").Value = 174.30500225754
$ws.Range("S10").Value = 0.03860714852883285
$ws.Range("T10").Value = 0.03860714852883285
